$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust column widths: col A -> 83, col C -> 49 (raw xlsx width units)
$ws.Columns.Item(1).ColumnWidth = 82.16666666666667
$ws.Columns.Item(3).ColumnWidth = 48.166666666666664

# Apply the shared style (wrap text + vertical center) used by the existing data rows
# to the new rows 4:5 before writing values, matching style index "1" in the target file.
$styleRange = $ws.Range("A4:G5")
$styleRange.WrapText = $true
$styleRange.VerticalAlignment = -4108

# Row 4
$ws.Range("A4").Value = "<20641191.1075855687472.JavaMail.evans@thyme>"
$ws.Range("B4").Value = "Tue, 17 Oct 2000 02:26:00 -0700 (PDT)"
$ws.Range("C4").Value = "phillip.allen@enron.com"
$ws.Range("D4").Value = "Phillip K Allen"
$ws.Range("E4").Value = "enron"
$ws.Range("F4").Value = "Re: High Speed Internet Access"
$ws.Range("G4").Value = "1. login:  pallen pw: ke9davis`n I don't think these are required by the ISP `n  2.  static IP address`n IP: 64.216.90.105`n Sub: 255.255.255.248`n gate: 64.216.90.110`n DNS: 151.164.1.8`n  3.  Company: 0413`n        RC:  105891"

# Row 5
$ws.Range("A5").Value = "`t<SI2PR06MB50915CD790AC930831A23BC0F77F2@SI2PR06MB5091.apcprd06.prod.outlook.com>"
$ws.Range("B5").Value = "Wed, 09 Oct 2024 10:52:18 +0000"
$ws.Range("C5").Value = "Krishnananda R <krishnananda.r@alphanimble.com>"
# D5 left blank (empty numeric cell in source)
$ws.Range("E5").Value = "alphanimble"
$ws.Range("F5").Value = "Introduction "
$ws.Range("G5").Value = "Hi, my name is R Krishnananda. You can use this mail address to communicate with me.`n(Please ignore mail from rkrishnananda2003@gmail.com)`nThank you`nR Krishnananda"

